# #327 Ajout des profils d'acces
# - Bump the "Date" metadata value.
# - Swap the two Mapping columns ("Mapping: RIM Mapping" and
#   "Mapping: Spécification métier vers l'extension ROR LevelRecourseORSAN"):
#   header text, row data and column widths all exchange places between
#   column AK (37) and column AL (38) on the "Elements" sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh generation Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- Elements sheet: swap Mapping columns AK <-> AL ---
$els = $wb.Worksheets.Item("Elements")

# Header row (row 1)
$els.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR LevelRecourseORSAN"
$els.Range("AL1").Value = "Mapping: RIM Mapping"

# Data rows (only rows 3, 5 and 6 carry values; rows 2 and 4 stay blank)
$els.Range("AK3").Value = ""
$els.Range("AL3").Value = "n/a"

$els.Range("AK5").Value = ""
$els.Range("AL5").Value = "N/A"

$els.Range("AK6").Value = "niveauRecoursORSAN"
$els.Range("AL6").Value = "N/A"

# Column widths: AK (37) and AL (38) trade widths too
$els.Columns("AK:AK").ColumnWidth = 76.5
$els.Columns("AL:AL").ColumnWidth = 24.166666666666664
